$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 11 (only 3 data rows remain: header + 2 users)
$ws.Range("A4:E11").EntireRow.Delete()

# Update remaining data rows with new values.
# B column holds numeric-looking IDs that must stay text, like the original
# inline-string cells. Entering them through a quoted formula and then
# pasting-as-values forces text storage without leaving the cell's number
# format / style changed.
$b2 = $ws.Range("B2")
$b2.Formula = "=""93528"""
$b2.Copy()
$b2.PasteSpecial(-4163)   # xlPasteValues

$ws.Range("C2").Value = "Ярослав Кузнецов"
$ws.Range("D2").Value = "?"

$b3 = $ws.Range("B3")
$b3.Formula = "=""49666"""
$b3.Copy()
$b3.PasteSpecial(-4163)   # xlPasteValues

$ws.Range("C3").Value = "Чернов Егор"
$ws.Range("D3").Value = "?"

$excel.CutCopyMode = $false
